$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells so numeric-looking strings
# (e.g. "0.0720", "4.70") are preserved exactly as text, matching the
# original inline-string cell contents instead of being coerced to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.339.43"
$ws.Range("E2").Value = "  -3.23%  "
$ws.Range("D3").Value = "2.244.98"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "237.01"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  -4.87%  "
$ws.Range("D7").Value = "69.56"
$ws.Range("E7").Value = "  -3.27%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -7.22%  "
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "58.94"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "36.87"
$ws.Range("E12").Value = "  +14.40%  "
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -5.70%  "
$ws.Range("D15").Value = "2.580.18"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").Value = "15.06"
$ws.Range("E16").Value = "  -6.18%  "
$ws.Range("D17").Value = "0.862"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "2.251.15"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").Value = "42.273.71"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("D22").Value = "73.42"
$ws.Range("E22").Value = "  -5.63%  "
$ws.Range("D23").Value = "236.46"
$ws.Range("E23").Value = "  -4.75%  "
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "10.02"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "171.13"
$ws.Range("D31").Value = "20.59"
$ws.Range("E31").Value = "  -6.97%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").Value = "  -5.16%  "
$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "3.71"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "22.41"
$ws.Range("E38").Value = "  +20.48%  "
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("D40").Value = "0.0275"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").Value = "5.93"
$ws.Range("E41").Value = "  -6.69%  "
$ws.Range("D42").Value = "65.06"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "9.36"
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  -14.47%  "
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("D46").Value = "0.192"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Value = "4.61"
$ws.Range("E47").Value = "  +13.32%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "10.15"
$ws.Range("E50").Value = "  +9.33%  "
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -2.82%  "
